# Apply the LinuxForHealth re-brand + version bump edit described by the diff.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-family-size"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the top-level "Extension" element; its Constraint(s) cell (column AI)
# loses its ele-1/ext-1 constraint text and becomes blank.
$elements.Range("AI2").Value = ""

# Row 5 ("Extension.url") carries the extension's canonical/fixed URL in
# column Q - same string as the Metadata sheet's URL, so it moves in lockstep.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-family-size"
